# NIT-9007719987.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker/period/value rows in the account-statement table (rows 16-33,
# columns B:G) are the same 18 (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora,
# Salario) records as before, just re-sorted so "Periodo Mora" (column E)
# runs in ascending order (ties keep their original relative order — i.e. a
# stable sort), instead of being grouped worker-by-worker with each worker's
# periods descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 33

# Read the existing table rows (columns B..G) into memory.
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()
    $g = $ws.Cells.Item($r, 7).Value()

    $row = @{
        B = $b
        C = $c
        D = $d
        E = $e
        F = $f
        G = $g
        Ord = [int]$e
    }
    $rows += $row
}

# Stable sort by Periodo Mora (column E) ascending.
$sorted = $rows | Sort-Object -Property Ord

# Write the reordered records back into the same range.
for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $firstRow + $i
    $row = $sorted[$i]
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
